# "Case 8: Merging within foreach" table had the merge-tag placeholder in
# its middle column changed from a bare "<<cellMerge>>" tag (no argument)
# followed by "<<[Country]>>", into "<<[Country]>>" followed by
# "<<cellMerge [Country]>>" -- i.e. the cellMerge instruction now names
# the field it merges on, and the plain [Country] field that used to
# trail it now comes first.
#
# Before: <<cellMerge>><<[Country]>>
# After : <<[Country]>><<cellMerge [Country]>>

$d = $word.ActiveDocument

# Find the right table by its distinctive "foreach [in clients]" text
# rather than assuming a fixed index.
$table = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    if ($candidate.Range.Text.Contains("foreach [in clients]")) {
        $table = $candidate
    }
}

$cell = $table.Cell(1, 2)
$rng = $cell.Range

$rng.Find.ClearFormatting()
$rng.Find.Execute("<<cellMerge>><<[Country]>>", $true, $false, $false, $false, $false, $true, 1, $false, `
    "<<[Country]>><<cellMerge [Country]>>", 2)
